$d = $word.ActiveDocument

function Replace-ParaByAnchor($anchor, $xml) {
    $rng = $d.Content
    $found = $rng.Find.Execute($anchor, $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
    if (-not $found) {
        Write-Output "NOT FOUND: $anchor"
        return
    }
    $p = $rng.Paragraphs(1).Range
    $p.InsertXML($xml)
}

# Title para - split into 3 spellchecked words
$xml = @'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" w14:paraId="1A58F12E" w14:textId="77777777" w:rsidR="00415EC3" w:rsidRPr="003F75FD" w:rsidRDefault="00C35233" w:rsidP="00C35233" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml"><w:pPr><w:pStyle w:val="Otsikko"/><w:jc w:val="both"/><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:b w:val="0"/><w:sz w:val="24"/><w:szCs w:val="24"/><w:lang w:val="en-US"/></w:rPr></w:pPr><w:proofErr w:type="spellStart"/><w:r w:rsidRPr="003F75FD"><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:b w:val="0"/><w:sz w:val="24"/><w:szCs w:val="24"/><w:lang w:val="en-US"/></w:rPr><w:t>Lappeenrannan</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:b w:val="0"/><w:sz w:val="24"/><w:szCs w:val="24"/><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:b w:val="0"/><w:sz w:val="24"/><w:szCs w:val="24"/><w:lang w:val="en-US"/></w:rPr><w:t>teknillinen</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:b w:val="0"/><w:sz w:val="24"/><w:szCs w:val="24"/><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:b w:val="0"/><w:sz w:val="24"/><w:szCs w:val="24"/><w:lang w:val="en-US"/></w:rPr><w:t>yliopisto</w:t></w:r><w:proofErr w:type="spellEnd"/></w:p>
'@
Replace-ParaByAnchor 'Lappeenrannan teknillinen yliopisto' $xml

# Date : 18.6.2024 -
$xml = @'
<w:p><w:pPr><w:rPr><w:lang w:val="en-US"/></w:rPr></w:pPr><w:proofErr w:type="gramStart"/><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>Date :</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r><w:r w:rsidR="00D91EFC"><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve">18.6.2024 - </w:t></w:r></w:p>
'@
Replace-ParaByAnchor 'Date : 18.6.2024 - ' $xml

# Activity : Self-learning
$xml = @'
<w:p><w:pPr><w:rPr><w:lang w:val="en-US"/></w:rPr></w:pPr><w:proofErr w:type="gramStart"/><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>Activity :</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r><w:r w:rsidR="00D91EFC"><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>Self-learning of Front-End development</w:t></w:r></w:p>
'@
Replace-ParaByAnchor 'Activity : Self-learning of Front-End development' $xml

# Date : ( date here)
$xml = @'
<w:p><w:pPr><w:rPr><w:lang w:val="en-US"/></w:rPr></w:pPr><w:proofErr w:type="gramStart"/><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>Date :</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve"> ( date here)</w:t></w:r></w:p>
'@
Replace-ParaByAnchor 'Date : ( date here)' $xml

# Activity : Video lecture
$xml = @'
<w:p><w:pPr><w:rPr><w:lang w:val="en-US"/></w:rPr></w:pPr><w:proofErr w:type="gramStart"/><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>Activity :</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve"> Video lecture or other activity details here</w:t></w:r></w:p>
'@
Replace-ParaByAnchor 'Activity : Video lecture or other activity details here' $xml

# mixins paragraph
$xml = @'
<w:p><w:pPr><w:pStyle w:val="Leipteksti"/><w:rPr><w:lang w:val="en-US"/></w:rPr></w:pPr><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve">Those advantages of SCSS that I hadn&#8217;t seen in the last task made themselves seen by allowing me to use variables and nest the styling of child elements inside the parents (which, while technically doesn&#8217;t let me DO anything new, does make the code much more readable), and yeah, now I see why SCSS is something I like to call &#8220;big good.&#8221; I also learned about </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>mixins</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>, which seem to basically act as functions for SCSS (and maybe they work in CSS as well? Not sure)</w:t></w:r></w:p>
'@
Replace-ParaByAnchor 'I also learned about mixins, which seem to basically act as functions for SCSS' $xml

# vh units / stackoverflow paragraph
$xml = @'
<w:p><w:pPr><w:pStyle w:val="Leipteksti"/><w:rPr><w:lang w:val="en-US"/></w:rPr></w:pPr><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve">I also found out about a few extra miscellaneous things I had not ran into yet, such as z-indexes and </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>vh</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve"> units. The material also had a very good explanation on the difference between rms and rems, which allowed me to understand it much better than previously</w:t></w:r><w:r w:rsidR="008D78D7"><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>.</w:t></w:r><w:r w:rsidR="008D78D7"><w:rPr><w:lang w:val="en-US"/></w:rPr><w:br/><w:t xml:space="preserve"> </w:t></w:r><w:r w:rsidR="008D78D7"><w:rPr><w:lang w:val="en-US"/></w:rPr><w:tab/><w:t xml:space="preserve">The one question that I was left with, however, was what is the difference between setting the height of something as 20% vs 20vh? Fortunately, the dynamic duo of Mr. Google and Mr. </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>Stackoverflow</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve"> were here to help and explain that </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>vh</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve"> can be used to always refer to the absolute size of the screen, while % is stuck referring to the height of its parent element. The same </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>stackoverflow</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve"> thread also confirmed what I suspected: There is also a </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>vw</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve">, a width equivalent of </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>vh</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve">. </w:t></w:r></w:p>
'@
Replace-ParaByAnchor 'I also found out about a few extra miscellaneous things' $xml

# Definitely makes paragraph
$xml = @'
<w:p><w:pPr><w:pStyle w:val="Leipteksti"/><w:rPr><w:lang w:val="en-US"/></w:rPr></w:pPr><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:br/></w:r><w:r w:rsidR="00CC1CAA"><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve">I learned about using JavaScript to affect the classes of HTML elements and about doing a rotate effect with SCSS on them. I had seen stuff like this online and always assumed they were done with some kind of animation file, but it turned out it can be done rather simply with the transform and translate functions. </w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>Definitely makes</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve"> creating these types of visual effects (and especially tinkering with them to get them to look just right) far easier than I had expected. </w:t></w:r></w:p>
'@
Replace-ParaByAnchor 'Definitely makes creating these types of visual effects' $xml

# _menu.scss paragraph
$xml = @'
<w:p><w:pPr><w:pStyle w:val="Leipteksti"/><w:rPr><w:lang w:val="en-US"/></w:rPr></w:pPr><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:tab/><w:t xml:space="preserve">Another interesting tidbit I came across was that in the video, the _</w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>menu.scss</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve"> file gains access to _</w:t></w:r><w:proofErr w:type="spellStart"/><w:proofErr w:type="gramStart"/><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>config.scss</w:t></w:r><w:proofErr w:type="spellEnd"/><w:proofErr w:type="gramEnd"/><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve"> from both of them being included in the </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>main.scss</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve"> file, without the _</w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>menu.scss</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve"> file needing to include the _</w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>config.scss</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve"> file.</w:t></w:r></w:p>
'@
Replace-ParaByAnchor 'Another interesting tidbit I came across was that in the video' $xml

# Stackoverflow Task 6 paragraph
$xml = @'
<w:p><w:pPr><w:pStyle w:val="Leipteksti"/><w:rPr><w:lang w:val="en-US"/></w:rPr></w:pPr><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:br/></w:r><w:r w:rsidR="00753A1C"><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve">For the first part where the work page was being made, there was once again, not a whole lot happening. The only new tidbit I learned was how to extend classes in (S)CSS. I had already known it was possible, technically, due to having seen it being mentioned in </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>Stackoverflow</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve"> during my research in a previous project, but simply had not needed to use it myself and therefore did not know the syntax for it.</w:t></w:r></w:p>
'@
Replace-ParaByAnchor 'For the first part where the work page was being made' $xml

# Project date merge
$xml = @'
<w:p><w:pPr><w:pStyle w:val="Leipteksti"/><w:rPr><w:lang w:val="en-US"/></w:rPr></w:pPr><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>24.6.2024 &#8211; 26.6.2024</w:t></w:r></w:p>
'@
Replace-ParaByAnchor '24.6.2024' $xml

# abandon its current slide + determineScroll
$xml = @'
<w:p><w:pPr><w:pStyle w:val="Leipteksti"/><w:rPr><w:lang w:val="en-US"/></w:rPr></w:pPr><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:tab/><w:t xml:space="preserve">In the case of the vertical scroll, the issues stemmed mainly from a problem of the mouse wheel being rotated enough to register twice causing the site to abandon its current slide and stop in the middle of the scrolling animation. While I had surmised that forcing a delay between the function calls for the scroll was the key to fix it, I continued to have the </w:t></w:r><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:lastRenderedPageBreak/><w:t xml:space="preserve">same issue even after doing that. It took me a while but after randomly moving the check from the </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>determineScroll</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve"> function to the </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>addEventListener</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>, it started working correctly.</w:t></w:r></w:p>
'@
Replace-ParaByAnchor 'abandon it' $xml

# Insert new 'mobile layout' paragraph + trailing empty paragraph at end of body
$xmlNew = @'
<w:p><w:pPr><w:pStyle w:val="Leipteksti"/><w:rPr><w:lang w:val="en-US"/></w:rPr></w:pPr><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:tab/><w:t xml:space="preserve">Another </w:t></w:r><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>problem</w:t></w:r><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve"> I </w:t></w:r><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>had</w:t></w:r><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve"> was the mobile layout. Not only was fitting some of the pages on my site into a small space rather difficult, but I also ran into an issue where I assume the top navigation bar of the browser messes with the positioning of some of my buttons. </w:t></w:r><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve">I found some </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>stackoverflow</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve"> threads talking about a similar issue and their suggestions about using </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>svh</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve"> instead of </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>vh</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve"> sadly did not work, so I tried simply moving the buttons further up and </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>and</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve"> down but for some reason it felt that no matter what I did, one of them was always out of view.</w:t></w:r><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve">In the end, I just </w:t></w:r><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>said,</w:t></w:r><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve"> &#8220;screw this&#8221; and removed the scroll buttons on small screens </w:t></w:r><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve">entirely </w:t></w:r><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>and made the mobile users rely on the menu for navigation.</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="Leipteksti"/><w:rPr><w:lang w:val="en-US"/></w:rPr></w:pPr></w:p>
'@
$lastPara = $d.Paragraphs($d.Paragraphs.Count)
$endRng = $d.Range($lastPara.Range.End, $lastPara.Range.End)
$endRng.InsertXML($xmlNew)

Write-Output 'DONE'